$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "latitude"
$ws.Range("E1").Value = "longitude"
$ws.Range("F1").Value = "coordinates"

$ws.Range("F2").Value = "-15.7600264,-47.8828874"
$ws.Range("F3").Value = "-15.7957539,-47.988938"
$ws.Range("F4").Value = "-15.7686722,-47.8888662"
$ws.Range("F5").Value = "-15.7634302,-47.8724814"

$ws.Range("F2").Select()
